# Applies the "LittleMedium" typography addition + translation sheet updates
# described by the commit diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Typography" (sheet1): add a new typography row (row 8) that
# defines the "LittleMedium" font entry.
# ---------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

$typo.Cells.Item(8, 2).Value = "LittleMedium"        # B8 - Typography Name
$typo.Cells.Item(8, 3).Value = "Roboto-Regular.ttf"  # C8 - Font
$typo.Cells.Item(8, 4).Value = 45                    # D8 - Size
$typo.Cells.Item(8, 5).Value = 4                     # E8 - Bpp
$typo.Cells.Item(8, 6).Value = "?"                   # F8 - Fallback Character
$typo.Cells.Item(8, 8).Value = "42-91"               # H8 - Wildcard Ranges
# G8 (Wildcard Characters) stays empty, same as the other rows in the table.
$typo.Range("B8:H8").Style = "Normal"

# ---------------------------------------------------------------------
# Sheet "Translation" (sheet2)
# ---------------------------------------------------------------------
$trans = $wb.Worksheets.Item("Translation")

# Rows 180-235: every entry that used the "Medium" typography now uses
# the new "LittleMedium" typography instead.
for ($r = 180; $r -le 235; $r++) {
    $trans.Cells.Item($r, 3).Value = "LittleMedium"
}

# Row 236: the GB wildcard placeholder text lost its leading "N".
$trans.Cells.Item(236, 5).Value = "<value>"

# Rows 240-243: new "Extra" typography rows.
$extraRows = 240..243
foreach ($r in $extraRows) {
    $trans.Cells.Item($r, 4).Value = "Left"
    $trans.Cells.Item($r, 5).Value = "N"
    $trans.Cells.Item($r, 6).Value = "LTR"
}
$trans.Cells.Item(240, 2).Value = "SingleUseId254"
$trans.Cells.Item(241, 2).Value = "SingleUseId255"
$trans.Cells.Item(242, 2).Value = "SingleUseId256"
$trans.Cells.Item(243, 2).Value = "SingleUseId257"
$trans.Range("C240:C243").Value = "Extra"

# Rows 244-257: new "LittleMedium" typography rows (debug screen texts).
$textIds = @{
    244 = "SingleUseId258"
    245 = "SingleUseId259"
    246 = "SingleUseId260"
    247 = "SingleUseId261"
    248 = "SingleUseId262"
    249 = "SingleUseId263"
    250 = "SingleUseId264"
    251 = "SingleUseId265"
    252 = "SingleUseId266"
    253 = "SingleUseId267"
    254 = "SingleUseId268"
    255 = "SingleUseId269"
    256 = "SingleUseId270"
    257 = "SingleUseId271"
}

# Alignment/GB alternate Center/<value> and Left/<wildcard> pairs, just
# like rows 180-235 do.
$alignments = @{
    244 = @("Center", "<value>")
    245 = @("Left", "ind_5")
    246 = @("Center", "<value>")
    247 = @("Left", "0.00")
    248 = @("Center", "<value>")
    249 = @("Left", "0.00")
    250 = @("Center", "<value>")
    251 = @("Left", "ind_5")
    252 = @("Center", "<value>")
    253 = @("Left", "0.00")
    254 = @("Center", "<value>")
    255 = @("Left", "ind_5")
    256 = @("Center", "<value>")
    257 = @("Left", "0.00")
}

foreach ($r in 244..257) {
    $trans.Cells.Item($r, 2).Value = $textIds[$r]
    $trans.Cells.Item($r, 3).Value = "LittleMedium"
    $trans.Cells.Item($r, 4).Value = $alignments[$r][0]
    $trans.Cells.Item($r, 5).Value = $alignments[$r][1]
    $trans.Cells.Item($r, 6).Value = "LTR"
}

# New cells should not carry an explicit style (matches the rest of the
# table, which relies on the column default formatting).
$trans.Range("B240:F257").Style = "Normal"
